$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected marking/total figures -------------------------------------
# C11 keeps its existing "incorrectStyle" (s=6) formatting and its text
# (inlineStr) type - only the displayed text changes from "-3" to "-1".
# Assigning the literal -1 straight to the cell would make Excel interpret
# it as a number, so stage the text in a scratch cell (forcing text via a
# leading apostrophe), copy it, and paste only the *value* into C11 - this
# carries over the text type without disturbing C11's own formatting.
$ws.Range("A13").Value = "'-1"
$ws.Range("A13").Copy()
$ws.Range("C11").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A13").Clear()
$excel.CutCopyMode = 0

$ws.Range("C12").Value = -7
$ws.Range("E12").Value = "83/140"

# --- Row-label formatting --------------------------------------------------
# A10 ("No."), A11 ("Marking") and A12 ("Total") pick up the same
# boxed/centered "mtitleStyle" formatting already used by the header row
# above them (row 9). Copy that formatting across so the existing cellXfs
# entry is reused instead of a new one being synthesized.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
